$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2224.35
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2224.35
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2224.35
$ws.Range("N32").Value = -2876.35
$ws.Range("M32").ClearContents()

$ws.Range("H111").Value = 4346.4
$ws.Range("I111").Value = 5175
$ws.Range("J111").Value = 1032
$ws.Range("K111").Value = 15525
$ws.Range("L111").Value = 3096
$ws.Range("M111").Value = -12458
$ws.Range("N111").Value = -9230

$ws.Range("H137").Value = 2327819
$ws.Range("I137").Value = 5001726
$ws.Range("J137").Value = 2682.6086
$ws.Range("K137").Value = 15005178
$ws.Range("L137").Value = 8047.825800000001
$ws.Range("M137").Value = -15002628
$ws.Range("N137").Value = -13147.8258

$ws.Range("H138").Value = 2036381
$ws.Range("I138").Value = 1748.909
$ws.Range("J138").Value = 2782412.8
$ws.Range("K138").Value = 5246.727000000001
$ws.Range("L138").Value = 8347238.399999999
$ws.Range("M138").Value = -106.7270000000008
$ws.Range("N138").Value = -8357518.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1699.875
$ws.Range("I2").Value = 1599.8334
$ws.Range("K2").Value = 1599.8334
$ws.Range("M2").Value = -1486.8334

$ws.Range("H32").Value = 4996132.5
$ws.Range("I32").Value = 5669475.5
$ws.Range("K32").Value = 5669475.5
$ws.Range("M32").Value = -5669188.5

$ws.Range("H61").Value = 55667484
$ws.Range("I61").Value = 66734190
$ws.Range("J61").Value = 333938
$ws.Range("K61").Value = 66734190
$ws.Range("L61").Value = 333938
$ws.Range("M61").Value = -66733978
$ws.Range("N61").Value = -334362

$ws.Range("H116").Value = 1699.875
$ws.Range("I116").Value = 1599.8334
$ws.Range("K116").Value = 1599.8334
$ws.Range("M116").Value = 694.1666

$ws.Range("H136").Value = 55667484
$ws.Range("I136").Value = 66734190
$ws.Range("J136").Value = 333938
$ws.Range("K136").Value = 200202570
$ws.Range("L136").Value = 1001814
$ws.Range("M136").Value = -200200020
$ws.Range("N136").Value = -1006914

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1699.875
$ws.Range("I3").Value = 1599.8334
$ws.Range("K3").Value = 1599.8334
$ws.Range("M3").Value = -1485.8334

$ws.Range("H64").Value = 966.6667
$ws.Range("J64").Value = 966.6667
$ws.Range("L64").Value = 966.6667
$ws.Range("N64").Value = -1416.6667

$ws.Range("H67").Value = 966.6667
$ws.Range("J67").Value = 966.6667
$ws.Range("L67").Value = 966.6667
$ws.Range("N67").Value = -2526.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 36.666668
$ws.Range("J7").Value = 163.33333
$ws.Range("K7").Value = 36.666668
$ws.Range("L7").Value = 163.33333
$ws.Range("M7").Value = 76.333332
$ws.Range("N7").Value = -389.33333

$ws.Range("H22").Value = 325.13333
$ws.Range("I22").Value = 288.63635
$ws.Range("J22").Value = 425.5
$ws.Range("K22").Value = 288.63635
$ws.Range("L22").Value = 425.5
$ws.Range("M22").Value = 61.36365000000001
$ws.Range("N22").Value = -1125.5

$ws.Range("H31").Value = 13529.8125
$ws.Range("I31").Value = 41017
$ws.Range("K31").Value = 41017
$ws.Range("M31").Value = -40722

$ws.Range("H34").Value = 13529.8125
$ws.Range("I34").Value = 41017
$ws.Range("K34").Value = 41017
$ws.Range("M34").Value = -40815

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2882.3
$ws.Range("I64").Value = 594.3333
$ws.Range("J64").Value = 3862.8572
$ws.Range("K64").Value = 1782.9999
$ws.Range("L64").Value = 11588.5716
$ws.Range("M64").Value = -1512.9999
$ws.Range("N64").Value = -12128.5716

$ws.Range("H67").Value = 2882.3
$ws.Range("I67").Value = 594.3333
$ws.Range("J67").Value = 3862.8572
$ws.Range("K67").Value = 1782.9999
$ws.Range("L67").Value = 11588.5716
$ws.Range("M67").Value = -846.9999
$ws.Range("N67").Value = -13460.5716

$ws.Range("H68").Value = 815.22076
$ws.Range("I68").Value = 473.92062
$ws.Range("J68").Value = 2351.0715
$ws.Range("K68").Value = 1421.76186
$ws.Range("L68").Value = 7053.2145
$ws.Range("M68").Value = -610.7618600000001
$ws.Range("N68").Value = -8675.2145

$ws.Range("H71").Value = 815.22076
$ws.Range("I71").Value = 473.92062
$ws.Range("J71").Value = 2351.0715
$ws.Range("K71").Value = 4265.28558
$ws.Range("L71").Value = 21159.6435
$ws.Range("M71").Value = -209.2855799999998
$ws.Range("N71").Value = -29271.6435

$ws.Range("H88").Value = 3513.3333
$ws.Range("J88").Value = 3513.3333
$ws.Range("L88").Value = 10539.9999
$ws.Range("N88").Value = -11395.9999

$ws.Range("H91").Value = 3513.3333
$ws.Range("J91").Value = 3513.3333
$ws.Range("L91").Value = 10539.9999
$ws.Range("N91").Value = -13503.9999

$ws.Range("H92").Value = 1136.25
$ws.Range("I92").Value = 775
$ws.Range("J92").Value = 1497.5
$ws.Range("K92").Value = 2325
$ws.Range("L92").Value = 4492.5
$ws.Range("M92").Value = -1077
$ws.Range("N92").Value = -6988.5

$ws.Range("H94").Value = 3068.75
$ws.Range("J94").Value = 3478.5715
$ws.Range("L94").Value = 10435.7145
$ws.Range("N94").Value = -11787.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 25000
$ws.Range("J108").Value = 25000
$ws.Range("L108").Value = 25000
$ws.Range("N108").Value = -32680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3416.5
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 3333
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 9999
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -14899

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3652.1538
$ws.Range("I122").Value = 3200
$ws.Range("J122").Value = 3934.75
$ws.Range("K122").Value = 9600
$ws.Range("L122").Value = 11804.25
$ws.Range("M122").Value = -7150
$ws.Range("N122").Value = -16704.25

$ws.Range("H136").Value = 63554.688
$ws.Range("I136").Value = 40954
$ws.Range("J136").Value = 144271.42
$ws.Range("K136").Value = 122862
$ws.Range("L136").Value = 432814.26
$ws.Range("M136").Value = -120312
$ws.Range("N136").Value = -437914.26
